$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh + symbol-list reorder (rows 42-43)
# matching the GitHub Actions "Updated symbol list" commit.
# Numeric-looking values are forced to stay Text (matching the source
# inlineStr cells) by briefly switching to a Text number format while the
# value is entered, then restoring the Normal style so no stray formatting
# is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.397'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05975'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.430'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.525'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8084'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9226'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1423'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07415'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03286'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03087'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09354'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.858'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001571'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04721'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005901'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17OneONE'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005875'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004886'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00006801'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.151'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3232'
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002340'
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03968'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006371'
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004100'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41CEJICEJIBestin24h'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1077'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009180'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '43LocalTradersLCT'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005058'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.7001'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002413'
$ws.Range("D48").Style = "Normal"
